$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in A1 (stored as serial date number)
$ws.Range("A1").Value = 45436

# Update price values in column D
$ws.Range("D33").Value = 277.464
$ws.Range("D34").Value = 307.393
$ws.Range("D35").Value = 355.55
$ws.Range("D39").Value = 379.403
$ws.Range("D40").Value = 570.456
